$d = $word.ActiveDocument

# Update the date/day header paragraph
$d.Content.Find.Execute("2023-10-21 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-10-22 Sunday", 2)

# Update the division problems table by addressing each cell directly
# (several new values duplicate other old values elsewhere in the table,
# so a blind global Find/Replace would corrupt already-replaced cells;
# direct cell addressing avoids that ambiguity entirely).
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)
$newValues = @(
    @("85÷4=", "53÷7=", "44÷6=", "73÷6=", "89÷9="),
    @("86÷7=", "58÷8=", "31÷4=", "57÷3=", "84÷6="),
    @("99÷5=", "93÷3=", "93÷9=", "98÷2=", "40÷8="),
    @("34÷8=", "89÷3=", "54÷7=", "37÷8=", "33÷7="),
    @("72÷8=", "14÷6=", "72÷3=", "43÷6=", "51÷9=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i][$c - 1]
    }
}
